$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "28.222.72"
$ws.Cells.Item(2, 5).Value = "  +1.46%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.793.39"
$ws.Cells.Item(3, 5).Value = "  +2.77%  "

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.005"
$ws.Cells.Item(4, 5).Value = "  +0.33%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "335.56"
$ws.Cells.Item(5, 5).Value = "  +0.63%  "

# Row 6
$ws.Cells.Item(6, 5).Value = "  +0.20%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4514"
$ws.Cells.Item(7, 5).Value = "  +16.25%  "

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3735"
$ws.Cells.Item(8, 5).Value = "  +10.74%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "45.08"
$ws.Cells.Item(9, 5).Value = "  -0.59%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.07559"
$ws.Cells.Item(10, 5).Value = "  +5.76%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "1.140"
$ws.Cells.Item(11, 5).Value = "  +3.75%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  +0.45%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "22.28"
$ws.Cells.Item(13, 5).Value = "  +1.87%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "6.279"
$ws.Cells.Item(14, 5).Value = "  +3.40%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "7.440"
$ws.Cells.Item(15, 5).Value = "  +7.25%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "1.792.25"
$ws.Cells.Item(16, 5).Value = "  +2.90%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  +3.64%  "

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.06724"
$ws.Cells.Item(18, 5).Value = "  +1.72%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "81.02"
$ws.Cells.Item(19, 5).Value = "  +2.57%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "1.002"
$ws.Cells.Item(20, 5).Value = "  +0.21%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "17.42"
$ws.Cells.Item(21, 5).Value = "  +4.10%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "6.363"

# Row 23
$ws.Cells.Item(23, 4).Value = "28.214.09"
$ws.Cells.Item(23, 5).Value = "  +1.55%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "11.75"
$ws.Cells.Item(24, 5).Value = "  +2.22%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.421"
$ws.Cells.Item(25, 5).Value = "  +1.29%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "20.52"
$ws.Cells.Item(26, 5).Value = "  +3.94%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "151.84"
$ws.Cells.Item(27, 5).Value = "  -1.23%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "2.353"
$ws.Cells.Item(28, 5).Value = "  +3.28%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "1.997.87"
$ws.Cells.Item(29, 5).Value = "  +2.97%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "133.21"
$ws.Cells.Item(30, 5).Value = "  +4.47%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  -3.46%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "4.038"
$ws.Cells.Item(32, 5).Value = "  -0.20%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.09401"
$ws.Cells.Item(33, 5).Value = "  +8.19%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  +0.59%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.2360"
$ws.Cells.Item(35, 5).Value = "  +13.42%  "

# Row 36
$ws.Cells.Item(36, 2).Value = "Aptos"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "12.04"
$ws.Cells.Item(36, 5).Value = "  +0.78%  "

# Row 37
$ws.Cells.Item(37, 2).Value = "Hedera"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.06310"
$ws.Cells.Item(37, 5).Value = "  +4.08%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.02325"
$ws.Cells.Item(38, 5).Value = "  +3.11%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "5.191"
$ws.Cells.Item(39, 5).Value = "  +1.93%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.6545"
$ws.Cells.Item(40, 5).Value = "  +2.17%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "WEMIXTOKEN"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "1.481"
$ws.Cells.Item(41, 5).Value = "  -1.69%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "FraxShare"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "8.313"
$ws.Cells.Item(42, 5).Value = "  +5.88%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "1.209"
$ws.Cells.Item(43, 5).Value = "  +1.68%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  +0.12%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "14.11"
$ws.Cells.Item(45, 5).Value = "  +4.12%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "3.833"
$ws.Cells.Item(46, 5).Value = "  +0.75%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.6065"
$ws.Cells.Item(47, 5).Value = "  +2.55%  "

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "129.52"
$ws.Cells.Item(48, 5).Value = "  +2.99%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "2.018"
$ws.Cells.Item(49, 5).Value = "  +2.51%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.07122"
$ws.Cells.Item(50, 5).Value = "  +3.10%  "

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.159"
$ws.Cells.Item(51, 5).Value = "  +1.45%  "
